$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F,H,I,K,L,N across rows 2-25
# (columns A, G, J, M are left unchanged)
$newData = @{
    2 = @{ "B" = 1.01932169267991; "C" = 0.04794862259097954; "D" = 0.09918881493162246; "E" = 0.0587433830858739; "F" = 1.921414720276005; "H" = 0.07973214163530429; "I" = 1.544883511641899; "K" = 0.814875740161483; "L" = 0.2467597928774126; "N" = 2.627067627999587 }
    3 = @{ "B" = 0.9749936956408192; "C" = 0.04172204689923831; "D" = 0.0993803596383529; "E" = 0.05837070214486673; "F" = 1.907890773193515; "H" = 0.07973214163530429; "I" = 1.542902292285319; "K" = 0.7674244833305863; "L" = 0.2392944416309604; "N" = 2.641760148719165 }
    4 = @{ "B" = 0.9483559185673585; "C" = 0.03789546538081368; "D" = 0.09949987815922956; "E" = 0.05816767181456584; "F" = 1.900532097282053; "H" = 0.07973214163530429; "I" = 1.542342171065279; "K" = 0.7387516762821349; "L" = 0.2348535731389063; "N" = 2.651500712145214 }
    5 = @{ "B" = 0.9376467304189191; "C" = 0.03633515587658565; "D" = 0.09954906736413527; "E" = 0.05809143077956946; "F" = 1.897770825965011; "H" = 0.07973214163530429; "I" = 1.54227888044948; "K" = 0.7271835219086995; "L" = 0.2330797952928663; "N" = 2.655650810055775 }
    6 = @{ "B" = 0.9358772951911476; "C" = 0.03607600813016631; "D" = 0.09955726458275738; "E" = 0.05807916357838216; "F" = 1.897326654564054; "H" = 0.07973214163530429; "I" = 1.542278332122656; "K" = 0.72526965870955; "L" = 0.2327874293391687; "N" = 2.656350843001313 }
    7 = @{ "B" = 0.9482108997035255; "C" = 0.03787442644021155; "D" = 0.09950053957521199; "E" = 0.05816661728980854; "F" = 1.90049389666072; "H" = 0.07973214163530429; "I" = 1.5423406496908; "K" = 0.7385951934766695; "L" = 0.2348295059661609; "N" = 2.651555950200219 }
    8 = @{ "B" = 1.003917122852442; "C" = 0.04580236202671983; "D" = 0.09925446586301412; "E" = 0.05860953331852059; "F" = 1.916555381448376; "H" = 0.07973214163530429; "I" = 1.544064106237592; "K" = 0.798418521257048; "L" = 0.2441560694355189; "N" = 2.631984238538251 }
    9 = @{ "B" = 1.117758773852415; "C" = 0.06132614393239066; "D" = 0.09878686389780533; "E" = 0.05968251765390775; "F" = 1.955563800163318; "H" = 0.07973214163530429; "I" = 1.55265736177823; "K" = 0.9194104010800856; "L" = 0.2635813068849018; "N" = 2.599317529893085 }
    10 = @{ "B" = 1.204216383256835; "C" = 0.07272466747211581; "D" = 0.09845214126896984; "E" = 0.06059527786076657; "F" = 1.988826758971058; "H" = 0.07973214163530429; "I" = 1.562159834447129; "K" = 1.010568793196938; "L" = 0.2785503289684641; "N" = 2.578808771767271 }
    11 = @{ "B" = 1.244163574549816; "C" = 0.07791027655727589; "D" = 0.09830172423477102; "E" = 0.06103751007249159; "F" = 2.004964152080134; "H" = 0.07973214163530429; "I" = 1.567177751167947; "K" = 1.052536748941179; "L" = 0.2855127060835656; "N" = 2.570238861047798 }
    12 = @{ "B" = 1.259379349812832; "C" = 0.07987409019091274; "D" = 0.09824502731026996; "E" = 0.0612088509519424; "F" = 2.011219958137659; "H" = 0.07973214163530429; "I" = 1.569178034118096; "K" = 1.068500995255988; "L" = 0.2881712189366823; "N" = 2.56710307045546 }
    13 = @{ "B" = 1.256098419745967; "C" = 0.07945114040072099; "D" = 0.09825722637208756; "E" = 0.06117177730639867; "F" = 2.009866208305823; "H" = 0.07973214163530429; "I" = 1.568742782876612; "K" = 1.065059610772664; "L" = 0.2875976812300678; "N" = 2.567773547834165 }
    14 = @{ "B" = 1.245413609659238; "C" = 0.0780718374799676; "D" = 0.09829705449718418; "E" = 0.06105152872752484; "F" = 2.005475914882638; "H" = 0.07973214163530429; "I" = 1.567340308556368; "K" = 1.053848695915605; "L" = 0.28573098223724; "N" = 2.569978682908555 }
    15 = @{ "B" = 1.23888039542436; "C" = 0.07722699481878692; "D" = 0.09832148449641842; "E" = 0.06097837770389347; "F" = 2.002805615278518; "H" = 0.07973214163530429; "I" = 1.566494293564702; "K" = 1.046991051157278; "L" = 0.2845904425615657; "N" = 2.571343651710734 }
    16 = @{ "B" = 1.201618150275749; "C" = 0.07238579018725488; "D" = 0.09846200831751162; "E" = 0.06056691989027385; "F" = 1.987792407402893; "H" = 0.07973214163530429; "I" = 1.561845902663237; "K" = 1.007836145608707; "L" = 0.2780984000529969; "N" = 2.579384139516222 }
    17 = @{ "B" = 1.178916967218811; "C" = 0.06941602052312135; "D" = 0.09854868655822635; "E" = 0.0603214173041664; "F" = 1.978840123445607; "H" = 0.07973214163530429; "I" = 1.559172419105124; "K" = 0.9839438650503212; "L" = 0.2741549201518865; "N" = 2.584511428736931 }
    18 = @{ "B" = 1.165917946412492; "C" = 0.06770791675342025; "D" = 0.09859871589943658; "E" = 0.06018275373365256; "F" = 1.973785670279568; "H" = 0.07973214163530429; "I" = 1.557700128545818; "K" = 0.9702486730897704; "L" = 0.2719011243271012; "N" = 2.587531996321019 }
    19 = @{ "B" = 1.161526683367867; "C" = 0.06712958463393193; "D" = 0.09861568501689177; "E" = 0.06013624162196152; "F" = 1.972090569508751; "H" = 0.07973214163530429; "I" = 1.557212869150085; "K" = 0.9656197897280094; "L" = 0.2711404984245718; "N" = 2.588566979993374 }
    20 = @{ "B" = 1.18132753466881; "C" = 0.06973215394862109; "D" = 0.098539441494502; "E" = 0.06034728833355985; "F" = 1.979783310004933; "H" = 0.07973214163530429; "I" = 1.559450243596942; "K" = 0.9864823727286591; "L" = 0.2745732208760359; "N" = 2.583958219908467 }
    21 = @{ "B" = 1.248549592554184; "C" = 0.0784769678219277; "D" = 0.09828534890608331; "E" = 0.06108674348335086; "F" = 2.006761515014006; "H" = 0.07973214163530429; "I" = 1.56774953127065; "K" = 1.057139661744657; "L" = 0.286278679379393; "N" = 2.569328009143121 }
    22 = @{ "B" = 1.292999712348205; "C" = 0.08419304116202397; "D" = 0.09812081477857859; "E" = 0.06159261578530462; "F" = 2.025238200155144; "H" = 0.07973214163530429; "I" = 1.573757161003272; "K" = 1.103737358089631; "L" = 0.2940571869794439; "N" = 2.560404373762594 }
    23 = @{ "B" = 1.269228611475626; "C" = 0.08114216134497099; "D" = 0.09820849077467919; "E" = 0.06132055710236983; "F" = 2.015299447268049; "H" = 0.07973214163530429; "I" = 1.570497333083679; "K" = 1.078828924985515; "L" = 0.2898939016379529; "N" = 2.565108631819243 }
    24 = @{ "B" = 1.180237554544249; "C" = 0.06958923235239922; "D" = 0.09854362057326149; "E" = 0.06033558431873587; "F" = 1.979356607754482; "H" = 0.07973214163530429; "I" = 1.559324437510384; "K" = 0.9853345863319305; "L" = 0.2743840654977561; "N" = 2.584208098778888 }
    25 = @{ "B" = 1.086467629072388; "C" = 0.05712826127383153; "D" = 0.09891179283161655; "E" = 0.05937038119142279; "F" = 1.94420432214568; "H" = 0.07973214163530429; "I" = 1.549773417740262; "K" = 0.886282705774903; "L" = 0.2582041647424944; "N" = 2.607542227327855 }
}

foreach ($row in $newData.Keys) {
    $rowVals = $newData[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}